# Apply "enhanced scraping" data-completeness updates to the
# Organizations worksheet: fill in previously-empty cells with newly
# scraped values, and widen the social/contact columns (G:L) to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column width updates (columns G through L, i.e. 7-12) ---
# Excel's ColumnWidth setter pads the stored OOXML <col width> by the
# default 5/6-character margin, so back the padding out here to land on
# the exact target stored widths (16, 43, 37, 36, 33, 43).
$pad = 5 / 6
$ws.Columns.Item(7).ColumnWidth  = 16 - $pad
$ws.Columns.Item(8).ColumnWidth  = 43 - $pad
$ws.Columns.Item(9).ColumnWidth  = 37 - $pad
$ws.Columns.Item(10).ColumnWidth = 36 - $pad
$ws.Columns.Item(11).ColumnWidth = 33 - $pad
$ws.Columns.Item(12).ColumnWidth = 43 - $pad

# --- Row 2: Student Life ---
$ws.Range("D2").Value = "https://beulah.edu/logos/studentlife_logo.png"
$ws.Range("E2").Value = "Student organization focused on general activities and community engagement. The Student Life welcomes all interested students to participate and make a positive impact."
$ws.Range("F2").Value = "studentlife@beulah.edu"
$ws.Range("G2").Value = "(555) 901-2345"
$ws.Range("J2").Value = "https://facebook.com/studentlife"
$ws.Range("K2").Value = "https://twitter.com/studentlife"

# --- Row 3: Student Groups ---
$ws.Range("I3").Value = "https://instagram.com/studentgroups"
$ws.Range("J3").Value = "https://facebook.com/studentgroups"

# --- Row 4: Life in Atlanta ---
$ws.Range("F4").Value = "lifeinatlanta@beulah.edu"
$ws.Range("G4").Value = "(555) 789-0123"
$ws.Range("J4").Value = "https://facebook.com/lifeinatlanta"
$ws.Range("L4").Value = "https://youtube.com/channel/lifeinatlanta"

# --- Row 5: Faculty Portal ---
$ws.Range("F5").Value = "facultyportal@beulah.edu"
$ws.Range("H5").Value = "https://linkedin.com/groups/facultyportal"
